$wb = $excel.ActiveWorkbook

# Sheet1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1749
$ws.Range("F3").Value = 10234
$ws.Range("F6").Value = 600
$ws.Range("F8").Value = 1639
$ws.Range("F9").Value = 188
$ws.Range("F10").Value = 396
$ws.Range("F12").Value = 216
$ws.Range("F14").Value = 490
$ws.Range("F17").Value = 33
$ws.Range("F18").Value = 15
$ws.Range("F19").Value = 13
$ws.Range("F20").Value = 101
$ws.Range("F21").Value = 365
$ws.Range("F23").Value = 335
$ws.Range("F24").Value = 105
$ws.Range("F25").Value = 1168
$ws.Range("F26").Value = 704
$ws.Range("F30").Value = 244
$ws.Range("F32").Value = 501
$ws.Range("F33").Value = 229
$ws.Range("F36").Value = 676
$ws.Range("F37").Value = 55807
$ws.Range("F38").Value = 761
$ws.Range("F40").Value = 1271
$ws.Range("F41").Value = 820
$ws.Range("F42").Value = 387
$ws.Range("F43").Value = 343
$ws.Range("F44").Value = 27
$ws.Range("F46").Value = 81

# Sheet2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 89
$ws.Range("F14").Value = 109
$ws.Range("F18").Value = 1106
$ws.Range("F20").Value = 1052
$ws.Range("F22").Value = 332
$ws.Range("F23").Value = 688
$ws.Range("F24").Value = 80
$ws.Range("F28").Value = 368
$ws.Range("F31").Value = 212
$ws.Range("F34").Value = 160
$ws.Range("F43").Value = 73

# Sheet3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 812
$ws.Range("F5").Value = 199
$ws.Range("F7").Value = 4134
$ws.Range("F10").Value = 334
$ws.Range("F11").Value = 203

# Sheet4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1749
$ws.Range("F3").Value = 812
$ws.Range("F5").Value = 10234
$ws.Range("F6").Value = 199
$ws.Range("F7").Value = 4134
$ws.Range("F9").Value = 334
$ws.Range("F10").Value = 334
$ws.Range("F11").Value = 600
$ws.Range("F12").Value = 1639
$ws.Range("F13").Value = 188
$ws.Range("F20").Value = 109
$ws.Range("F22").Value = 15
$ws.Range("F23").Value = 101
$ws.Range("F24").Value = 1106
$ws.Range("F25").Value = 365
$ws.Range("F26").Value = 335
$ws.Range("F28").Value = 704
$ws.Range("F29").Value = 80
$ws.Range("F32").Value = 244
$ws.Range("F33").Value = 368
$ws.Range("F34").Value = 501
$ws.Range("F39").Value = 212
$ws.Range("F40").Value = 761
$ws.Range("F42").Value = 820
$ws.Range("F43").Value = 387
$ws.Range("F44").Value = 343
$ws.Range("F45").Value = 27
$ws.Range("F49").Value = 73
